# Update res_bus vm_pu results for the 380 kV slack-voltage case (commit "case with 380 kV done").
# Column B is the slack-bus voltage setpoint (changed from 1.05 to 1.02 p.u.); columns C:F and I:N are
# the recalculated per-bus voltage magnitudes for rows 2-25. Columns G (constant 1) and H (blank) are
# untouched, matching the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> ordered [B, C, D, E, F, I, J, K, L, M, N] values
$newValues = @{
    2 = @(1.02, 1.033366029227542, 1.026513451475762, 1.041296853725589, 1.049333393264743, 1.028140619206563, 1.038491151223592, 1.029335878222186, 1.044076768678686, 1.052090743175996, 1.016359438818662)
    3 = @(1.02, 1.035456500123056, 1.026832767717716, 1.043210675826196, 1.051452051245787, 1.028119731977348, 1.040218527726439, 1.029464047931006, 1.045798295184488, 1.054018255724803, 1.016972109842518)
    4 = @(1.02, 1.036804338530668, 1.027038990663314, 1.044444645201334, 1.052818661432289, 1.028104546995959, 1.041331346126495, 1.029545830079758, 1.046907440178718, 1.055260824523872, 1.017365833892181)
    5 = @(1.02, 1.037369841693726, 1.027125589626533, 1.044962378131683, 1.053392182731227, 1.028097762610924, 1.041798024285772, 1.029579933844558, 1.047372600427122, 1.055782111355798, 1.017530713907021)
    6 = @(1.02, 1.03746472686782, 1.027140124207022, 1.045049248248012, 1.053488421571141, 1.028096599963332, 1.04187631483863, 1.029585643700014, 1.047450637651736, 1.055869574615349, 1.01755836060787)
    7 = @(1.02, 1.036811899204402, 1.027040148187642, 1.0444515671777, 1.052826328749198, 1.028104457917983, 1.0413375863911, 1.029546286868134, 1.046913660057774, 1.055267794227486, 1.017368039536576)
    8 = @(1.02, 1.034073532951003, 1.026621446480098, 1.041944565930772, 1.050050311958026, 1.02813390540539, 1.039075957746021, 1.029379431388346, 1.04465957474345, 1.052743136263629, 1.016567061135038)
    9 = @(1.02, 1.029209801128964, 1.025880711243438, 1.037492025308604, 1.045124365226372, 1.02817304963526, 1.035051958522594, 1.029076656696341, 1.040649743577527, 1.048257462646569, 1.015134458890293)
    10 = @(1.02, 1.025939617291363, 1.025385058414673, 1.034498556500015, 1.041815550827719, 1.028190631074268, 1.032341682945349, 1.02886901662911, 1.037949528146233, 1.045240498997322, 1.014164613834184)
    11 = @(1.02, 1.024516597594767, 1.025170031100778, 1.033196024245224, 1.040376494231811, 1.028196236062541, 1.031161200832587, 1.028777753015332, 1.036773553826061, 1.043927448861385, 1.013741033680938)
    12 = @(1.02, 1.023986934918851, 1.025090101753861, 1.032711220915079, 1.039840979805725, 1.028198017265132, 1.030721647574989, 1.028743651841878, 1.036335698255753, 1.043438686032643, 1.013583141053325)
    13 = @(1.02, 1.024100599274606, 1.025107249484422, 1.032815257965506, 1.039955894605354, 1.028197648785466, 1.030815982143403, 1.028750975765279, 1.036429667602273, 1.043543574808476, 1.013617034894785)
    14 = @(1.02, 1.024472837959735, 1.025163425298077, 1.033155970530062, 1.040332248738053, 1.028196389426981, 1.031124889223288, 1.028774938309027, 1.036737382070019, 1.043887068933161, 1.013727993633358)
    15 = @(1.02, 1.024702040976964, 1.025198029374265, 1.033365763308621, 1.040564001381583, 1.028195573671256, 1.031315074403884, 1.028789675730796, 1.036926835544606, 1.044098568512557, 1.013796284965319)
    16 = @(1.02, 1.02603390762909, 1.025399320714895, 1.034584864946816, 1.041910920141611, 1.028190216864643, 1.032419879352862, 1.028875045105079, 1.038027428624988, 1.04532749814416, 1.014192648055195)
    17 = @(1.019999999999999, 1.02686745018941, 1.025525478374841, 1.035347855050772, 1.042754090270336, 1.028186319765488, 1.03311102144153, 1.02892823368536, 1.038715969682559, 1.046096560209594, 1.014440296640642)
    18 = @(1.02, 1.027352967761147, 1.025599024633332, 1.035792284708298, 1.043245289960371, 1.02818385285959, 1.033513487843186, 1.02895912692504, 1.039116933547572, 1.046544498368184, 1.014584396272636)
    19 = @(1.02, 1.027518403541067, 1.025624095254975, 1.035943721312831, 1.043412674438442, 1.028182978809425, 1.033650606598615, 1.028969638492465, 1.039253542438254, 1.04669712568265, 1.014633471543508)
    20 = @(1.02, 1.026778088910315, 1.025511946906639, 1.035266056715309, 1.042663689154085, 1.028186757925705, 1.033036937450378, 1.028922540567341, 1.038642163243583, 1.046014113847933, 1.014413762527972)
    21 = @(1.02, 1.024363253327552, 1.025146884520153, 1.03305566662043, 1.040221449297817, 1.028196768571165, 1.03103395356437, 1.028767887494854, 1.036646796958579, 1.043785947385488, 1.01369533448756)
    22 = @(1.02, 1.022838626403462, 1.024917017289605, 1.031660192195374, 1.038680200795731, 1.028201323127522, 1.029768393883549, 1.028669484195791, 1.035386161784353, 1.042378990933125, 1.013240407459075)
    23 = @(1.02, 1.023647471882663, 1.025038905493158, 1.032400512166615, 1.039497799546493, 1.028199073251875, 1.030439889504659, 1.028721759700692, 1.036055033862706, 1.043125427002243, 1.013481881958444)
    24 = @(1.02, 1.026818469506129, 1.025518061316702, 1.035303019744479, 1.042704539402101, 1.028186560538906, 1.033070414871011, 1.028925113446604, 1.038675515205249, 1.046051369794693, 1.01442575323124)
    25 = @(1.02, 1.030471933579147, 1.026072541693807, 1.03864741647204, 1.046402085227728, 1.028164434695646, 1.03609700809811, 1.029155958733625, 1.041691023941352, 1.049421663745579, 1.015507385343153)
}

$columns = @("B", "C", "D", "E", "F", "I", "J", "K", "L", "M", "N")

foreach ($row in $newValues.Keys) {
    $rowValues = $newValues[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $rowValues[$i]
    }
}
